# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# All Price (D) / Volume(1h) (E) cells hold plain text in the source data
# (e.g. "66.601.55", "  -4.36%  "), and rows 31/32 swap their USDe /
# NEARProtocol content. Leading apostrophes force numeric-looking Price
# values to stay text instead of being parsed into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.601.55"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "3.347.58"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'573.65"
$ws.Range("D6").Value = "'181.72"
$ws.Range("E6").Value = "  -5.11%  "
$ws.Range("E7").Value = "  +2.89%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.128"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "'0.404"
$ws.Range("E11").Value = "  -3.06%  "
$ws.Range("D12").Value = "3.936.87"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'26.89"
$ws.Range("E14").Value = "  -5.74%  "
$ws.Range("D15").Value = "66.725.04"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").Value = "3.357.58"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'435.86"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("D19").Value = "'13.61"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "'5.68"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "'7.60"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").Value = "'73.49"
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'0.519"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'9.06"
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "'1.95"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("D30").Value = "'22.86"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'5.29"
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("D33").Value = "'6.79"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").Value = "'1.22"
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("D35").Value = "'160.92"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -5.30%  "
$ws.Range("D37").Value = "'27.96"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  -7.62%  "
$ws.Range("D39").Value = "2.814.06"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").Value = "'0.799"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "'4.44"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").Value = "'6.21"
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("D43").Value = "'40.42"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "'0.0669"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "'24.30"
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("D46").Value = "'2.34"
$ws.Range("E46").Value = "  -6.97%  "
$ws.Range("D47").Value = "'326.20"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").Value = "'0.0273"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D50").Value = "'0.976"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").Value = "'6.15"
$ws.Range("E51").Value = "  -2.45%  "

# Reset style to Normal to drop the quotePrefix marker added by the leading apostrophe,
# while keeping the values stored as text (matches source data which are plain strings).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
